# Applies a cyclic rotation of species-record data between rows 2, 3 and 4:
#   new row2 <- old row4
#   new row3 <- old row2
#   new row4 <- old row3
# for columns A, B, E, F, G, H, Q, R (the remaining columns are identical
# across these three rows, so no visible change happens to them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes what used to be row 4.
$ws.Range("A2").Value = 112181853
$ws.Range("B2").Value = 78242
$ws.Range("E2").Value = 6453
$ws.Range("F2").Value = "Vedskivlav"
$ws.Range("G2").Value = "Hertelidea botryosa"
$ws.Range("H2").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q2").Value = 431106
$ws.Range("R2").Value = 6811802

# Row 3 becomes what used to be row 2.
$ws.Range("A3").Value = 112182534
$ws.Range("B3").Value = 77402
$ws.Range("E3").Value = 6446
$ws.Range("F3").Value = "Kolflarnlav"
$ws.Range("G3").Value = "Carbonicola anthracophila"
$ws.Range("H3").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 431104
$ws.Range("R3").Value = 6811805

# Row 4 becomes what used to be row 3.
$ws.Range("A4").Value = 112181898
$ws.Range("B4").Value = 78216
$ws.Range("E4").Value = 229821
$ws.Range("F4").Value = "Vedflamlav"
$ws.Range("G4").Value = "Ramboldia elabens"
$ws.Range("H4").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q4").Value = 431104
$ws.Range("R4").Value = 6811804
